$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Cxcl12 -> Itgav LR-pair results (natmi re-run per Dr Hou's advice):
# the sending/target cluster set now also includes the "M2" cluster, so the
# 3x3 block (ECs/FAPs/sCs) grows to a 4x4 block (ECs/FAPs/M2/sCs) and all
# specificity-derived figures are recomputed across the new cluster set.
$data = @(
    @("ECs", "Cxcl12", "Itgav", "ECs", 3, 1, 198.977211, 596.9316329999999, 0.440791350614085, 0.4407913506140851, 3, 1, 13.441269, 40.323807, 0.0897308213348123, 0.08973082133481232, 2674.506217920759, 24070.55596128683, 0.03955256992788307, 0.03955256992788309),
    @("ECs", "Cxcl12", "Itgav", "FAPs", 3, 1, 198.977211, 596.9316329999999, 0.440791350614085, 0.4407913506140851, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 10886.34864382888, 97977.13779445995, 0.1609953505096373, 0.1609953505096373),
    @("ECs", "Cxcl12", "Itgav", "M2", 3, 1, 198.977211, 596.9316329999999, 0.440791350614085, 0.4407913506140851, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 12670.29375234021, 114032.6437710619, 0.1873776461196113, 0.1873776461196114),
    @("ECs", "Cxcl12", "Itgav", "sCs", 3, 1, 198.977211, 596.9316329999999, 0.440791350614085, 0.4407913506140851, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 3574.732778005992, 32172.59500205392, 0.05286578405695329, 0.0528657840569533),
    @("FAPs", "Cxcl12", "Itgav", "ECs", 3, 1, 196.7746736666667, 590.324021, 0.4359121013721307, 0.4359121013721308, 3, 1, 13.441269, 40.323807, 0.0897308213348123, 0.08973082133481232, 2644.901321140883, 23804.11189026795, 0.03911475088590525, 0.03911475088590526),
    @("FAPs", "Cxcl12", "Itgav", "FAPs", 3, 1, 196.7746736666667, 590.324021, 0.4359121013721307, 0.4359121013721308, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 10765.84444542741, 96892.6000088467, 0.1592132455730546, 0.1592132455730546),
    @("FAPs", "Cxcl12", "Itgav", "M2", 3, 1, 196.7746736666667, 590.324021, 0.4359121013721307, 0.4359121013721308, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 12530.04254028644, 112770.3828625779, 0.1853035077851939, 0.185303507785194),
    @("FAPs", "Cxcl12", "Itgav", "sCs", 3, 1, 196.7746736666667, 590.324021, 0.4359121013721307, 0.4359121013721308, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 3535.163008379215, 31816.46707541294, 0.0522805971279769, 0.05228059712797692),
    @("M2", "Cxcl12", "Itgav", "ECs", 2, 0.6666666666666666, 0.158424, 0.475272, 0.0003509544061791369, 0.0003509544061791369, 3, 1, 13.441269, 40.323807, 0.0897308213348123, 0.08973082133481232, 2.129419600056, 19.164776400504, 0.00003149142711752528, 0.00003149142711752528),
    @("M2", "Cxcl12", "Itgav", "FAPs", 2, 0.6666666666666666, 0.158424, 0.475272, 0.0003509544061791369, 0.0003509544061791369, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 8.66762022084, 78.00858198756, 0.0001281831586690537, 0.0001281831586690538),
    @("M2", "Cxcl12", "Itgav", "M2", 2, 0.6666666666666666, 0.158424, 0.475272, 0.0003509544061791369, 0.0003509544061791369, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 10.08798247464, 90.79184227176, 0.0001491885229452398, 0.0001491885229452398),
    @("M2", "Cxcl12", "Itgav", "sCs", 2, 0.6666666666666666, 0.158424, 0.475272, 0.0003509544061791369, 0.0003509544061791369, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 2.846172497728, 25.615552479552, 0.00004209129744731807, 0.00004209129744731808),
    @("sCs", "Cxcl12", "Itgav", "ECs", 3, 1, 55.49875533333334, 166.496266, 0.122945593607605, 0.122945593607605, 3, 1, 13.441269, 40.323807, 0.0897308213348123, 0.08973082133481232, 745.9736996005181, 6713.763296404663, 0.01103200909390644, 0.01103200909390644),
    @("sCs", "Cxcl12", "Itgav", "FAPs", 3, 1, 55.49875533333334, 166.496266, 0.122945593607605, 0.122945593607605, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 3036.422094876104, 27327.79885388493, 0.04490484876551318, 0.04490484876551318),
    @("sCs", "Cxcl12", "Itgav", "M2", 3, 1, 55.49875533333334, 166.496266, 0.122945593607605, 0.122945593607605, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 3534.000348223754, 31806.00313401378, 0.05226340285234086, 0.05226340285234087),
    @("sCs", "Cxcl12", "Itgav", "sCs", 3, 1, 55.49875533333334, 166.496266, 0.122945593607605, 0.122945593607605, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 997.0650348928731, 8973.585314035858, 0.01474533289584447, 0.01474533289584447)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
